$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.36
$wsSummary.Range("B4").Value = 0.36
$wsSummary.Range("B5").Value = 0.12
$wsSummary.Range("B6").Value = 62
$wsSummary.Range("B7").Value = 20
$wsSummary.Range("B9").Value = 32.26

$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.36
$wsStatus.Range("D4").Value = 62
$wsStatus.Range("E4").Value = 0.36
$wsStatus.Range("F4").Value = 0.36
$wsStatus.Range("G4").Value = 32.26

$newRow = @(62, "2026-02-17", "15:44:20", "MarketMaking", "DOWN", 0.9, 0.95, "CLOSED", 5.5556, 0.05, 100.36, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsAllTrades.Range("B63").NumberFormat = "@"
for ($i = 0; $i -lt $newRow.Length; $i++) {
    $wsAllTrades.Cells.Item(63, $i + 1).Value = $newRow[$i]
}

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
$wsMarketMaking.Range("B63").NumberFormat = "@"
for ($i = 0; $i -lt $newRow.Length; $i++) {
    $wsMarketMaking.Cells.Item(63, $i + 1).Value = $newRow[$i]
}
